# Apply the committed changes: replace the TODAY() formula in B2 with a
# static date value, update the date value in B3, update the active
# selection on the sheet, and update the window position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# B2 previously held =TODAY() (a volatile formula); replace with the
# static date value 2002-12-03 (serial 37593).
$ws.Range("B2").Value = 37593

# B3 date value updated from 2023-12-10 (45270) to 2004-12-22 (38343).
$ws.Range("B3").Value = 38343

# Update the active cell / selection to B3.
$ws.Range("B3").Select()

# Update workbook window position.
$excel.ActiveWindow.Left = 3570
$excel.ActiveWindow.Top = 2145

$wb.Save()
